$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (column D) and "Volume(1h)" (column E) values pulled
# from coinranking.com. Rows 10 and 11 (Polygon / Dogecoin) also swap order.
#
# A literal leading apostrophe (the standard Excel "store as text" prefix)
# is used for price strings that would otherwise be auto-parsed as numbers
# (e.g. "2.310"), so they keep their original text representation/formatting
# exactly like the source data, which stores every cell as text.

$updates = @(
    @{ Cell = 'D2'; Value = '23.371.25' }
    @{ Cell = 'E2'; Value = '  -0.75%  ' }
    @{ Cell = 'D3'; Value = '1.625.93' }
    @{ Cell = 'E3'; Value = '  -0.99%  ' }
    @{ Cell = 'D4'; Value = '''0.9993' }
    @{ Cell = 'E4'; Value = '  +0.05%  ' }
    @{ Cell = 'D5'; Value = '''0.9998' }
    @{ Cell = 'E5'; Value = '  +0.14%  ' }
    @{ Cell = 'D6'; Value = '''304.38' }
    @{ Cell = 'E6'; Value = '  -1.13%  ' }
    @{ Cell = 'D7'; Value = '''0.3782' }
    @{ Cell = 'E7'; Value = '  -0.04%  ' }
    @{ Cell = 'D8'; Value = '''51.92' }
    @{ Cell = 'E8'; Value = '  -2.23%  ' }
    @{ Cell = 'D9'; Value = '''0.3623' }
    @{ Cell = 'E9'; Value = '  -1.73%  ' }
    @{ Cell = 'B10'; Value = 'Dogecoin' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge' }
    @{ Cell = 'D10'; Value = '''0.08097' }
    @{ Cell = 'E10'; Value = '  -1.19%  ' }
    @{ Cell = 'B11'; Value = 'Polygon' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic' }
    @{ Cell = 'D11'; Value = '''1.226' }
    @{ Cell = 'E11'; Value = '  -4.55%  ' }
    @{ Cell = 'D12'; Value = '''0.9995' }
    @{ Cell = 'E12'; Value = '  +0.03%  ' }
    @{ Cell = 'D13'; Value = '''22.71' }
    @{ Cell = 'E13'; Value = '  -2.45%  ' }
    @{ Cell = 'D14'; Value = '''6.551' }
    @{ Cell = 'E14'; Value = '  -1.89%  ' }
    @{ Cell = 'E15'; Value = '  -3.02%  ' }
    @{ Cell = 'D16'; Value = '''7.218' }
    @{ Cell = 'E16'; Value = '  -3.51%  ' }
    @{ Cell = 'D17'; Value = '1.624.58' }
    @{ Cell = 'E17'; Value = '  -1.00%  ' }
    @{ Cell = 'D18'; Value = '''93.51' }
    @{ Cell = 'E18'; Value = '  -1.68%  ' }
    @{ Cell = 'D19'; Value = '''0.06907' }
    @{ Cell = 'E19'; Value = '  -0.73%  ' }
    @{ Cell = 'D20'; Value = '''17.86' }
    @{ Cell = 'E20'; Value = '  -3.32%  ' }
    @{ Cell = 'D21'; Value = '''0.9997' }
    @{ Cell = 'E21'; Value = '  +0.14%  ' }
    @{ Cell = 'E22'; Value = '  -2.86%  ' }
    @{ Cell = 'D23'; Value = '23.364.50' }
    @{ Cell = 'E23'; Value = '  -0.76%  ' }
    @{ Cell = 'E24'; Value = '  -2.46%  ' }
    @{ Cell = 'D25'; Value = '''3.245' }
    @{ Cell = 'E25'; Value = '  +3.57%  ' }
    @{ Cell = 'E26'; Value = '  +0.65%  ' }
    @{ Cell = 'D27'; Value = '''21.07' }
    @{ Cell = 'E27'; Value = '  -1.94%  ' }
    @{ Cell = 'D28'; Value = '''149.85' }
    @{ Cell = 'E28'; Value = '  -1.38%  ' }
    @{ Cell = 'D29'; Value = '''5.291' }
    @{ Cell = 'E29'; Value = '  -0.80%  ' }
    @{ Cell = 'D30'; Value = '''134.09' }
    @{ Cell = 'E30'; Value = '  -1.85%  ' }
    @{ Cell = 'D31'; Value = '''2.310' }
    @{ Cell = 'E31'; Value = '  -4.86%  ' }
    @{ Cell = 'D32'; Value = '1.804.80' }
    @{ Cell = 'E32'; Value = '  -0.62%  ' }
    @{ Cell = 'D33'; Value = '''6.789' }
    @{ Cell = 'E33'; Value = '  -1.07%  ' }
    @{ Cell = 'D34'; Value = '''11.02' }
    @{ Cell = 'E34'; Value = '  +5.45%  ' }
    @{ Cell = 'D35'; Value = '''0.9513' }
    @{ Cell = 'E35'; Value = '  -2.78%  ' }
    @{ Cell = 'E36'; Value = '  -1.97%  ' }
    @{ Cell = 'D37'; Value = '''0.2516' }
    @{ Cell = 'E37'; Value = '  -1.22%  ' }
    @{ Cell = 'D38'; Value = '''0.08827' }
    @{ Cell = 'E38'; Value = '  -0.48%  ' }
    @{ Cell = 'D39'; Value = '''6.084' }
    @{ Cell = 'E39'; Value = '  -2.47%  ' }
    @{ Cell = 'E40'; Value = '  -4.88%  ' }
    @{ Cell = 'E41'; Value = '  -2.85%  ' }
    @{ Cell = 'D42'; Value = '''0.7056' }
    @{ Cell = 'E42'; Value = '  -1.81%  ' }
    @{ Cell = 'D43'; Value = '''16.22' }
    @{ Cell = 'E43'; Value = '  -0.22%  ' }
    @{ Cell = 'D44'; Value = '''12.29' }
    @{ Cell = 'E44'; Value = '  -3.19%  ' }
    @{ Cell = 'D45'; Value = '''0.6452' }
    @{ Cell = 'E45'; Value = '  -3.11%  ' }
    @{ Cell = 'D46'; Value = '''2.317' }
    @{ Cell = 'E46'; Value = '  -2.25%  ' }
    @{ Cell = 'D47'; Value = '''0.9988' }
    @{ Cell = 'E47'; Value = '  +0.10%  ' }
    @{ Cell = 'D48'; Value = '''3.993' }
    @{ Cell = 'E48'; Value = '  -1.28%  ' }
    @{ Cell = 'D49'; Value = '''0.07987' }
    @{ Cell = 'E49'; Value = '  -1.11%  ' }
    @{ Cell = 'D50'; Value = '''1.199' }
    @{ Cell = 'E50'; Value = '  -1.42%  ' }
    @{ Cell = 'D51'; Value = '''125.69' }
    @{ Cell = 'E51'; Value = '  -4.70%  ' }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
